# Applies the revisions described in the commit message:
# - Risk-of-bias (domain) judgements updated for several arms (rows 4-15)
# - Updated randomized_n / response_e / response_n / dropout_any_e / dropout_any_n
#   and overall PANSS summary stats for the two NCT phase III trials (rows 6-11)
# - Swapped death/serious-event counts for row 13
# - Removed two stray overall_baseline values (rows 18-19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that should be stored as TEXT even though it looks
# like a number (mirrors the source file, where these columns use
# inline/shared strings rather than numeric cells). We temporarily force a
# text number format so Excel does not auto-convert the value to a number,
# then restore the cell's original style so no stray formatting is left
# behind.
function Set-TextValue {
    param($addr, $value)
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = "$value"
    $r.Style = $origStyle
}

# Helper: write a plain (non-numeric-looking) text value.
function Set-PlainText {
    param($addr, $value)
    $ws.Range($addr).Value = $value
}

# Helper: write a numeric value.
function Set-NumValue {
    param($addr, $value)
    $ws.Range($addr).Value = $value
}

# ---- Row 4 ----
Set-PlainText "AE4" "Low"

# ---- Row 5 ----
Set-PlainText "AE5" "Low"

# ---- Row 6 (ulotaront 50mg/d (adults)) ----
Set-NumValue  "J6"  144
Set-TextValue "N6"  "22"
Set-TextValue "O6"  "144"
Set-TextValue "R6"  "34"
Set-TextValue "S6"  "144"
Set-PlainText "AB6" "Some concerns"
Set-PlainText "AD6" "Low"
Set-PlainText "AE6" "Low"
Set-NumValue  "AG6" 142
Set-NumValue  "AH6" 102
Set-NumValue  "AJ6" 19.07
Set-PlainText "AM6" "yes"

# ---- Row 7 (ulotaront 75mg/d (adults)) ----
Set-TextValue "N7"  "29"
Set-TextValue "R7"  "27"
Set-TextValue "S7"  "145"
Set-PlainText "AB7" "Some concerns"
Set-PlainText "AD7" "Low"
Set-PlainText "AE7" "Low"
Set-NumValue  "AH7" 102
Set-NumValue  "AJ7" 19.27
Set-PlainText "AM7" "yes"

# ---- Row 8 (placebo (adults)) ----
Set-NumValue  "J8"  146
Set-TextValue "N8"  "26"
Set-TextValue "O8"  "146"
Set-TextValue "R8"  "27"
Set-TextValue "S8"  "146"
Set-PlainText "AB8" "Some concerns"
Set-PlainText "AD8" "Low"
Set-PlainText "AE8" "Low"
Set-NumValue  "AH8" 102
Set-NumValue  "AJ8" 18.06
Set-PlainText "AM8" "yes"

# ---- Row 9 (ulotaront 75mg/d) ----
Set-NumValue  "J9"  155
Set-TextValue "N9"  "23"
Set-TextValue "O9"  "155"
Set-TextValue "R9"  "34"
Set-TextValue "S9"  "155"
Set-PlainText "AB9" "Some concerns"
Set-PlainText "AD9" "Low"
Set-PlainText "AE9" "Low"
Set-NumValue  "AG9" 153
Set-NumValue  "AH9" 101
Set-NumValue  "AJ9" 18.55
Set-PlainText "AM9" "yes"

# ---- Row 10 (ulotaront 100mg/d) ----
Set-TextValue "N10"  "27"
Set-TextValue "R10"  "38"
Set-TextValue "S10"  "154"
Set-PlainText "AB10" "Some concerns"
Set-PlainText "AD10" "Low"
Set-PlainText "AE10" "Low"
Set-NumValue  "AG10" 152
Set-NumValue  "AH10" 100
Set-NumValue  "AJ10" 18.49
Set-PlainText "AM10" "yes"

# ---- Row 11 (placebo) ----
Set-NumValue  "J11"  155
Set-TextValue "N11"  "21"
Set-TextValue "O11"  "155"
Set-TextValue "R11"  "27"
Set-TextValue "S11"  "155"
Set-PlainText "AB11" "Some concerns"
Set-PlainText "AD11" "Low"
Set-PlainText "AE11" "Low"
Set-NumValue  "AG11" 155
Set-NumValue  "AH11" 100
Set-NumValue  "AJ11" 18.67
Set-PlainText "AM11" "yes"

# ---- Row 12 ----
Set-PlainText "AE12" "Low"

# ---- Row 13 (death_e / serious_e swap) ----
Set-TextValue "V13" "0"
Set-TextValue "X13" "2"
Set-PlainText "AE13" "Low"

# ---- Row 14 ----
Set-PlainText "AE14" "Low"

# ---- Row 15 ----
Set-PlainText "AE15" "Low"

# ---- Row 18: drop stray overall_baseline value ----
$ws.Range("AH18").ClearContents()

# ---- Row 19: drop stray overall_baseline value ----
$ws.Range("AH19").ClearContents()
